$wb = $excel.ActiveWorkbook

# --- ALC row 15 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 633.5454999999999
$ws.Range("I15").Value = 633.5454999999999
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1900.6365
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -1731.6365

# --- ALC row 88 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1649.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1649.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 1649.5
$ws.Range("N88").Value = -2461.5
$ws.Range("M88").ClearContents()

# --- ALC row 91 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1649.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1649.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1649.5
$ws.Range("N91").Value = -4457.5
$ws.Range("M91").ClearContents()

# --- ALC row 96 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1999
$ws.Range("I96").Value = 1999
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 5997
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -4624
$ws.Range("N96").ClearContents()

# --- ALC row 100 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5000
$ws.Range("I100").Value = 5000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 5000
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -4459
$ws.Range("N100").ClearContents()

# --- ALC row 116 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5500
$ws.Range("I116").Value = 5500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 5500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -2058

# --- ARM row 37 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 23250
$ws.Range("I37").Value = 22200
$ws.Range("J37").Value = 25000
$ws.Range("K37").Value = 22200
$ws.Range("L37").Value = 25000
$ws.Range("M37").Value = -21927

# --- BSM row 7 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 189
$ws.Range("I7").Value = 48.333332
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 48.333332
$ws.Range("L7").Value = 400
$ws.Range("M7").Value = 64.666668
$ws.Range("N7").Value = -626

# --- BSM row 22 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 380
$ws.Range("I22").Value = 380
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 380
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -207

# --- BSM row 107 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1071.7778
$ws.Range("I107").Value = 665.5714
$ws.Range("J107").Value = 2493.5
$ws.Range("K107").Value = 665.5714
$ws.Range("L107").Value = 2493.5
$ws.Range("M107").Value = 1254.4286

# --- BSM row 134 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9930.823
$ws.Range("I134").Value = 8676.5
$ws.Range("J134").Value = 30000
$ws.Range("K134").Value = 26029.5
$ws.Range("L134").Value = 90000
$ws.Range("M134").Value = -23494.5

# --- CRP row 31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5392.7334
$ws.Range("I31").Value = 4066.5
$ws.Range("J31").Value = 6908.4287
$ws.Range("K31").Value = 4066.5
$ws.Range("L31").Value = 6908.4287
$ws.Range("M31").Value = -3771.5
$ws.Range("N31").Value = -7498.4287

# --- CRP row 34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5392.7334
$ws.Range("I34").Value = 4066.5
$ws.Range("J34").Value = 6908.4287
$ws.Range("K34").Value = 4066.5
$ws.Range("L34").Value = 6908.4287
$ws.Range("M34").Value = -3864.5
$ws.Range("N34").Value = -7312.4287

# --- CRP row 58 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 874
$ws.Range("I58").Value = 874
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 874
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -671
$ws.Range("N58").ClearContents()

# --- CRP row 136 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 874
$ws.Range("I136").Value = 874
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2622
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -72
$ws.Range("N136").ClearContents()

# --- CUL row 38 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 491.66666
$ws.Range("I38").Value = 491.66666
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1474.99998
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -1127.99998

# --- GSM row 11 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2705291
$ws.Range("I11").Value = 527.5
$ws.Range("J11").Value = 2855555.5
$ws.Range("K11").Value = 527.5
$ws.Range("L11").Value = 2855555.5
$ws.Range("M11").Value = -388.5
$ws.Range("N11").Value = -2855833.5

# --- GSM row 33 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 24999.666
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 24999.666
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 24999.666
$ws.Range("N33").Value = -25503.666

# --- GSM row 70 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2214.25
$ws.Range("I70").Value = 2214.25
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2214.25
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -1944.25
$ws.Range("N70").ClearContents()

# --- GSM row 73 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 2214.25
$ws.Range("I73").Value = 2214.25
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2214.25
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1278.25
$ws.Range("N73").ClearContents()

# --- GSM row 122 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6946384
$ws.Range("I122").Value = 7354668.5
$ws.Range("J122").Value = 5555
$ws.Range("K122").Value = 22064005.5
$ws.Range("L122").Value = 16665
$ws.Range("M122").Value = -22061555.5

# --- LTW row 22 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2546.8572
$ws.Range("I22").Value = 2798.625
$ws.Range("J22").Value = 2211.1667
$ws.Range("K22").Value = 2798.625
$ws.Range("L22").Value = 2211.1667
$ws.Range("M22").Value = -2503.625
$ws.Range("N22").Value = -2801.1667

# --- LTW row 27 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2546.8572
$ws.Range("I27").Value = 2798.625
$ws.Range("J27").Value = 2211.1667
$ws.Range("K27").Value = 2798.625
$ws.Range("L27").Value = 2211.1667
$ws.Range("M27").Value = -2691.625
$ws.Range("N27").Value = -2425.1667

# --- LTW row 122 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3302.4
$ws.Range("I122").Value = 3253
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 9759
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -7309

# --- WVR row 62 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4869.8
$ws.Range("I62").Value = 4855.3335
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4855.3335
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4231.3335

# --- WVR row 65 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4869.8
$ws.Range("I65").Value = 4855.3335
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 24276.6675
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -21156.6675

# --- WVR row 118 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 225000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 225000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 225000
$ws.Range("N118").Value = -228314
$ws.Range("M118").ClearContents()

# --- WVR row 122 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 987.25
$ws.Range("I122").Value = 984.7143
$ws.Range("J122").Value = 1005
$ws.Range("K122").Value = 2954.1429
$ws.Range("L122").Value = 3015
$ws.Range("M122").Value = -504.1428999999998

# --- WVR row 136 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 865
$ws.Range("I136").Value = 968.8889
$ws.Range("J136").Value = 397.5
$ws.Range("K136").Value = 2906.6667
$ws.Range("L136").Value = 1192.5
$ws.Range("M136").Value = -356.6667000000002
$ws.Range("N136").Value = -6292.5
